# Adds "death data" (Baseline deaths23, Population23, Death rate) columns
# H, I, J to the NCD+LRI block (rows 2-14) of Sheet1 in GEMMcoefficients.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# Order matters: Excel appends new shared-string entries in the order
# cells are first written, and the target file has the shared strings
# appended as: Population23, Baseline deaths23, Death rate.
$ws.Range("I1").Value = "Population23"
$ws.Range("H1").Value = "Baseline deaths23"
$ws.Range("J1").Value = "Death rate"

# --- Row 2 (0-25, built from published sub age-band totals) ----------
$ws.Range("H2").Formula = "=286+46+47+83+104"
$ws.Range("I2").Formula = "=176300+193300+194300+166300+157300"
$ws.Range("J2").Formula = "=H2/I2"

# --- Rows 3-13: plain baseline-deaths / population values ------------
$ws.Range("H3").Value = 116
$ws.Range("I3").Value = 181100

$ws.Range("H4").Value = 181
$ws.Range("I4").Value = 234200

$ws.Range("H5").Value = 242
$ws.Range("I5").Value = 224800

$ws.Range("H6").Value = 347
$ws.Range("I6").Value = 202600

$ws.Range("H7").Value = 500
$ws.Range("I7").Value = 170400

$ws.Range("H8").Value = 736
$ws.Range("I8").Value = 163000

$ws.Range("H9").Value = 1304
$ws.Range("I9").Value = 183600

$ws.Range("H10").Value = 2403
$ws.Range("I10").Value = 229300

$ws.Range("H11").Value = 2952
$ws.Range("I11").Value = 182500

$ws.Range("H12").Value = 3252
$ws.Range("I12").Value = 130700

$ws.Range("H13").Value = 2400
$ws.Range("I13").Value = 57300

# --- Row 14 (80+, combination of two published bands) -----------------
$ws.Range("H14").Formula = "=3656+5658"
$ws.Range("I14").Value = 90700

# --- Death rate column for rows 3-14, entered as one shared formula ---
$ws.Range("J3:J14").Formula = "=H3/I3"

# Restore the selection to match the saved workbook view.
$ws.Range("K14").Select()
